$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '34.867.82'
$ws.Range("E2").Value = '  +1.00%  '
$ws.Range("D3").Value = '1.809.29'
$ws.Range("E3").Value = '  +0.40%  '
$ws.Range("E4").Value = '  +0.31%  '
$ws.Range("D5").Value = '''231.87'
$ws.Range("E5").Value = '  +3.24%  '
$ws.Range("E6").Value = '  +0.70%  '
$ws.Range("E7").Value = '  +0.39%  '
$ws.Range("D8").Value = '''40.16'
$ws.Range("E8").Value = '  -5.12%  '
$ws.Range("E9").Value = '  +7.27%  '
$ws.Range("D10").Value = '''0.0683'
$ws.Range("E10").Value = '  +2.58%  '
$ws.Range("D11").Value = '''0.0995'
$ws.Range("E11").Value = '  +0.01%  '
$ws.Range("D12").Value = '2.071.16'
$ws.Range("E12").Value = '  +0.45%  '
$ws.Range("D13").Value = '1.815.27'
$ws.Range("E13").Value = '  +0.28%  '
$ws.Range("D14").Value = '''11.06'
$ws.Range("D15").Value = '''4.68'
$ws.Range("E15").Value = '  +6.12%  '
$ws.Range("E16").Value = '  +4.39%  '
$ws.Range("D17").Value = '34.830.45'
$ws.Range("E17").Value = '  +1.06%  '
$ws.Range("D18").Value = '''68.96'
$ws.Range("E18").Value = '  +2.51%  '
$ws.Range("D19").Value = '0.0₃0782'
$ws.Range("E19").Value = '  +1.78%  '
$ws.Range("D20").Value = '''236.66'
$ws.Range("E20").Value = '  -1.48%  '
$ws.Range("D21").Value = '''11.75'
$ws.Range("E21").Value = '  +5.35%  '
$ws.Range("E22").Value = '  +6.17%  '
$ws.Range("E23").Value = '  +0.33%  '
$ws.Range("E24").Value = '  +5.74%  '
$ws.Range("D25").Value = '''172.80'
$ws.Range("E25").Value = '  +1.06%  '
$ws.Range("D26").Value = '''7.90'
$ws.Range("E26").Value = '  +3.09%  '
$ws.Range("E27").Value = '  -0.04%  '
$ws.Range("E28").Value = '  -0.60%  '
$ws.Range("E29").Value = '  +31.46%  '
$ws.Range("D30").Value = '''1.00'
$ws.Range("E30").Value = '  +0.22%  '
$ws.Range("D31").Value = '3.339.01'
$ws.Range("E31").Value = '  +37.43%  '
$ws.Range("E32").Value = '  +7.16%  '
$ws.Range("E33").Value = '  +2.02%  '
$ws.Range("E34").Value = '  +1.89%  '
$ws.Range("E35").Value = '  -0.19%  '
$ws.Range("D36").Value = '''1.15'
$ws.Range("E36").Value = '  +8.51%  '
$ws.Range("D37").Value = '''92.99'
$ws.Range("E37").Value = '  +5.99%  '
$ws.Range("E38").Value = '  +4.18%  '
$ws.Range("E39").Value = '  +2.20%  '
$ws.Range("D40").Value = '1.304.03'
$ws.Range("E40").Value = '  -1.09%  '
$ws.Range("E41").Value = '  +4.41%  '
$ws.Range("D42").Value = '''0.985'
$ws.Range("E42").Value = '  +5.15%  '
$ws.Range("D43").Value = '''14.76'
$ws.Range("E43").Value = '  -0.48%  '
$ws.Range("D44").Value = '''2.33'
$ws.Range("E44").Value = '  -0.37%  '
$ws.Range("E45").Value = '  +0.47%  '
$ws.Range("E46").Value = '  -1.37%  '
$ws.Range("D47").Value = '''6.26'
$ws.Range("E47").Value = '  +7.95%  '
$ws.Range("E48").Value = '  -1.26%  '
$ws.Range("D49").Value = '1.986.07'
$ws.Range("E49").Value = '  +1.11%  '
$ws.Range("D50").Value = '''1.00'
$ws.Range("E50").Value = '  +0.30%  '
$ws.Range("E51").Value = '  +5.67%  '
